# Applies the edits described by the commit:
# "Defined client class and executed unit test"
#
# The workbook is a unit-test plan spreadsheet (BankAccount class).
# This script:
#   1. Fills in the developer/student name (C3).
#   2. Fills in the Preconditions / Method Inputs / Expected Result
#      columns (E, F, G) for test-case rows 7-22, which were previously
#      left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Student / developer name -------------------------------------------
$ws.Range("C3").Value2 = "Jashanpreet Kaur Jattana"

# --- Test case data (rows 7-22 => columns E=Preconditions, F=Method Inputs, G=Expected Result) ---

# Row 7: __init__ - Attributes are set to input values.
$ws.Range("E7").Value2 = "none"
$ws.Range("F7").Value2 = "account_number=12345, client_number=67890, balance=500.0"
$ws.Range("G7").Value2 = "account_number=12345, client_number=67890, balance=500.00"

# Row 8: __init__ - Balance attribute set to 0 when non-numeric balance argument.
$ws.Range("E8").Value2 = "none"
$ws.Range("F8").Value2 = 'account_number=12345, client_number=67890, balance="abc"'
$ws.Range("G8").Value2 = "Balance should be set to 0.0"

# Row 9: __init__ - ValueError when non-numeric account number
$ws.Range("E9").Value2 = "none"
$ws.Range("F9").Value2 = 'account_number="abc", client_number=67890, balance=500.0'
$ws.Range("G9").Value2 = "ValueError because account_number must be an integer"

# Row 10: __init__ - ValueError when non-numeric client number
$ws.Range("E10").Value2 = "none"
$ws.Range("F10").Value2 = 'account_number=12345, client_number="abc", balance=500.0'
$ws.Range("G10").Value2 = "ValueError because client_number must be an integer"

# Row 11: account_number (getter) - returns account number attribute
$ws.Range("E11").Value2 = "BankAccount object with account_number=12345"
$ws.Range("F11").Value2 = "none"
$ws.Range("G11").Value2 = 12345

# Row 12: client_number (getter) - returns client number attribute
$ws.Range("E12").Value2 = "BankAccount object with client_number=67890"
$ws.Range("F12").Value2 = "none"
$ws.Range("G12").Value2 = 67890

# Row 13: balance (getter) - returns balance attribute
$ws.Range("E13").Value2 = "BankAccount object with balance=500.0"
$ws.Range("F13").Value2 = "none"
$ws.Range("G13").Value2 = 500

# Row 14: update_balance - correctly updates balance attribute when positive amount is received.
$ws.Range("E14").Value2 = "BankAccount object with balance=500.0"
$ws.Range("F14").Value2 = "amount=100.0"
$ws.Range("G14").Value2 = "Updated balance should be 600.0"

# Row 15: update_balance - correctly updates balance attribute when negative amount is received.
$ws.Range("E15").Value2 = "BankAccount object with balance=500.0"
$ws.Range("F15").Value2 = "amount=-50.0"
$ws.Range("G15").Value2 = "Updated balance should be 450.0"

# Row 16: update_balance - Balance attribute value remains unchanged when amount is non-numeric
$ws.Range("E16").Value2 = "BankAccount object with balance=500.0"
$ws.Range("F16").Value2 = 'amount="abc"'
$ws.Range("G16").Value2 = "Balance should remain 500.0"

# Row 17: __str__ - returns string in expected format.
$ws.Range("E17").Value2 = "BankAccount object with balance=500.0"
$ws.Range("F17").Value2 = "amount=200.0"
$ws.Range("G17").Value2 = "Updated balance should be 700.0"

# Row 18: deposit - BankAccount object's balance is updated correctly when a valid amount is provided.
$ws.Range("E18").Value2 = "BankAccount object with balance=500.0"
$ws.Range("F18").Value2 = "amount=-50.0"
$ws.Range("G18").Value2 = "ValueError because deposit amount must be positive"

# Row 19: deposit - ValueError when negative amount is provided.
$ws.Range("E19").Value2 = "BankAccount object with balance=500.0"
$ws.Range("F19").Value2 = "amount=200.0"
$ws.Range("G19").Value2 = "Updated balance should be 300.0"

# Row 20: withdraw - ValueError when amount exceeds balance.
$ws.Range("E20").Value2 = "BankAccount object with balance=500.0"
$ws.Range("F20").Value2 = "amount=-50.0"
$ws.Range("G20").Value2 = "ValueError because withdrawl amount must be positive"

# Row 21: withdraw
$ws.Range("E21").Value2 = "BankAccount object with balance=500.0"
$ws.Range("F21").Value2 = "amount=600.0"
$ws.Range("G21").Value2 = "ValueError because withdrawl amount must not exceed the current balance"

# Row 22: __str__
$ws.Range("E22").Value2 = "BankAccount object with balance=500.0"
$ws.Range("F22").Value2 = "none"
$ws.Range("G22").Value2 = '"Account Number: 12345 Balance: 500"'

# --- Update the selection to match the final state of the authored file ---
$ws.Range("G22").Select()
